# Feature: add arrows (arrow_n). Fixed bugs, removed unnecessary code.
#
# On the "meta" sheet, row 5 (A5) used to be an empty, styled placeholder
# cell. We now fill that row with a new "style" / "default" property pair,
# and push the empty styled placeholder down to a new row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("meta")

# Copy the formatting of the old placeholder cell (A5) down to the new
# placeholder row (A6) before we overwrite A5 with real content.
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the new "style" / "default" row.
$ws.Range("A5").Value = "style"
$ws.Range("B5").Value = "default"
